# Applies the edit described by the diff:
#  - Within several groups of fixtures sharing the same kickoff timestamp,
#    the match-detail columns (F:V) were re-shuffled between rows (the
#    Indice/pais/torneio/temporada/data_partida columns A:E stay put).
#  - Three brand-new fixture rows (129-131) were appended at the bottom.
#  - The sheet's <dimension> grows from A1:V128 to A1:V131 (handled
#    automatically by Excel once the new cells are populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Pair($r1, $r2) {
    $rng1 = $ws.Range("F$r1`:V$r1")
    $rng2 = $ws.Range("F$r2`:V$r2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# Simple pairwise swaps of F:V content (same kickoff time, order flipped)
Swap-Pair 16 17
Swap-Pair 33 34
Swap-Pair 44 47
Swap-Pair 63 64
Swap-Pair 98 99
Swap-Pair 103 106
Swap-Pair 107 108
Swap-Pair 113 114

# Cyclic rotations of F:V content among >2 rows sharing a kickoff time
# Group: 28 -> 31 -> 30 -> 29 -> 28  (new(r) gets old content of the row listed)
$v28 = $ws.Range("F28:V28").Value()
$v29 = $ws.Range("F29:V29").Value()
$v30 = $ws.Range("F30:V30").Value()
$v31 = $ws.Range("F31:V31").Value()
$ws.Range("F28:V28").Value = $v31
$ws.Range("F29:V29").Value = $v28
$ws.Range("F30:V30").Value = $v29
$ws.Range("F31:V31").Value = $v30

# Group: 35 -> 36 -> 37 -> 35
$v35 = $ws.Range("F35:V35").Value()
$v36 = $ws.Range("F36:V36").Value()
$v37 = $ws.Range("F37:V37").Value()
$ws.Range("F35:V35").Value = $v36
$ws.Range("F36:V36").Value = $v37
$ws.Range("F37:V37").Value = $v35

# Group: 90 -> 91 -> 92 -> 93 -> 90
$v90 = $ws.Range("F90:V90").Value()
$v91 = $ws.Range("F91:V91").Value()
$v92 = $ws.Range("F92:V92").Value()
$v93 = $ws.Range("F93:V93").Value()
$ws.Range("F90:V90").Value = $v91
$ws.Range("F91:V91").Value = $v92
$ws.Range("F92:V92").Value = $v93
$ws.Range("F93:V93").Value = $v90

# Append three new fixture rows at the bottom, copying formatting from the
# last existing row (128) so number formats / styles (bold Indice column,
# date-formatted data_partida column) carry over correctly.
$ws.Range("A128:V128").Copy()
$ws.Range("A129:V131").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A129").Value = 128
$ws.Range("B129").Value = "italy"
$ws.Range("C129").Value = "serie-c-group-c"
$ws.Range("D129").Value = "2023-2024"
$ws.Range("E129").Value = 45243.85416666666
$ws.Range("F129").Value = "Picerno"
$ws.Range("G129").Value = 3
$ws.Range("H129").Value = "Potenza"
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = 1.88
$ws.Range("K129").Value = "09/11/2023 09:13"
$ws.Range("L129").Value = 2.34
$ws.Range("M129").Value = "13/11/2023 20:28"
$ws.Range("N129").Value = 3.14
$ws.Range("O129").Value = "09/11/2023 09:13"
$ws.Range("P129").Value = 2.91
$ws.Range("Q129").Value = "13/11/2023 20:20"
$ws.Range("R129").Value = 3.96
$ws.Range("S129").Value = "09/11/2023 09:13"
$ws.Range("T129").Value = 3.49
$ws.Range("U129").Value = "13/11/2023 20:28"
$ws.Range("V129").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/picerno-potenza/8f1EHXx8/"

$ws.Range("A130").Value = 129
$ws.Range("B130").Value = "italy"
$ws.Range("C130").Value = "serie-c-group-c"
$ws.Range("D130").Value = "2023-2024"
$ws.Range("E130").Value = 45243.86458333334
$ws.Range("F130").Value = "Benevento"
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = "Giugliano"
$ws.Range("I130").Value = 2
$ws.Range("J130").Value = 1.44
$ws.Range("K130").Value = "09/11/2023 09:13"
$ws.Range("L130").Value = 1.35
$ws.Range("M130").Value = "13/11/2023 20:41"
$ws.Range("N130").Value = 3.97
$ws.Range("O130").Value = "09/11/2023 09:13"
$ws.Range("P130").Value = 4.57
$ws.Range("Q130").Value = "13/11/2023 20:41"
$ws.Range("R130").Value = 6.42
$ws.Range("S130").Value = "09/11/2023 09:13"
$ws.Range("T130").Value = 10.56
$ws.Range("U130").Value = "13/11/2023 20:43"
$ws.Range("V130").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/benevento-giugliano/vF9cLFwq/"

$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "italy"
$ws.Range("C131").Value = "serie-c-group-c"
$ws.Range("D131").Value = "2023-2024"
$ws.Range("E131").Value = 45243.86458333334
$ws.Range("F131").Value = "Juve Stabia"
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = "Foggia"
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 2.02
$ws.Range("K131").Value = "09/11/2023 09:13"
$ws.Range("L131").Value = 1.77
$ws.Range("M131").Value = "13/11/2023 20:44"
$ws.Range("N131").Value = 2.99
$ws.Range("O131").Value = "09/11/2023 09:13"
$ws.Range("P131").Value = 3.49
$ws.Range("Q131").Value = "13/11/2023 20:44"
$ws.Range("R131").Value = 3.66
$ws.Range("S131").Value = "09/11/2023 09:13"
$ws.Range("T131").Value = 4.88
$ws.Range("U131").Value = "13/11/2023 20:44"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/juve-stabia-foggia/U32AIDN1/"

Write-Host "Done. UsedRange=$($ws.UsedRange.Address())"
